{"js": "// Replace each old math-fact / date string with its updated value.\n// Each (old, new) pair corresponds to one <w:t> run in the document\n// (the header date paragraph + the 100 table-cell expressions).\nconst pairs = [\n  [\"2023-11-23 Thursday\", \"2023-11-24 Friday\"],\n  [\"64+2=66\", \"15+64=79\"],\n  [\"11+47=58\", \"48-5=43\"],\n  [\"22+53=75\", \"40-27=13\"],\n  [\"91-32=59\", \"94-27=67\"],\n  [\"32+20=52\", \"94-42=52\"],\n  [\"69-11=58\", \"15+32=47\"],\n  [\"62-61=1\", \"17+79=96\"],\n  [\"18+34=52\", \"43+14=57\"],\n  [\"4+26=30\", \"85+1=86\"],\n  [\"43+11=54\", \"59-59=0\"],\n  [\"55+39=94\", \"78-13=65\"],\n  [\"10-10=0\", \"53-4=49\"],\n  [\"50+0=50\", \"30+7=37\"],\n  [\"79+10=89\", \"21+32=53\"],\n  [\"28-24=4\", \"65+9=74\"],\n  [\"56-27=29\", \"49+13=62\"],\n  [\"54-50=4\", \"91-29=62\"],\n  [\"84-24=60\", \"63-45=18\"],\n  [\"38+5=43\", \"6+7=13\"],\n  [\"22+22=44\", \"57+21=78\"],\n  [\"59+23=82\", \"80-44=36\"],\n  [\"64-12=52\", \"39-22=17\"],\n  [\"8+71=79\", \"22+38=60\"],\n  [\"28-15=13\", \"23+8=31\"],\n  [\"34+21=55\", \"4+75=79\"],\n  [\"3+61=64\", \"67-22=45\"],\n  [\"14+4=18\", \"94-73=21\"],\n  [\"87-35=52\", \"63+15=78\"],\n  [\"28+40=68\", \"71+1=72\"],\n  [\"63+32=95\", \"10+13=23\"],\n  [\"88-26=62\", \"78-6=72\"],\n  [\"79-46=33\", \"71-17=54\"],\n  [\"75-75=0\", \"7+12=19\"],\n  [\"99-6=93\", \"57-33=24\"],\n  [\"60-12=48\", \"54-8=46\"],\n  [\"27+36=63\", \"64+28=92\"],\n  [\"28+24=52\", \"18+52=70\"],\n  [\"5+76=81\", \"25+63=88\"],\n  [\"45+31=76\", \"68-3=65\"],\n  [\"27+58=85\", \"94-29=65\"],\n  [\"61+7=68\", \"20-1=19\"],\n  [\"98-65=33\", \"61+34=95\"],\n  [\"44-35=9\", \"34+59=93\"],\n  [\"51-2=49\", \"70-52=18\"],\n  [\"4+34=38\", \"65-46=19\"],\n  [\"94+3=97\", \"74-34=40\"],\n  [\"23+73=96\", \"69-30=39\"],\n  [\"48-27=21\", \"27+56=83\"],\n  [\"25+51=76\", \"29+19=48\"],\n  [\"83-8=75\", \"39+42=81\"],\n  [\"94-81=13\", \"67-10=57\"],\n  [\"9+42=51\", \"2+26=28\"],\n  [\"49-38=11\", \"58+20=78\"],\n  [\"81-23=58\", \"46+16=62\"],\n  [\"9+28=37\", \"93-45=48\"],\n  [\"34-12=22\", \"88-47=41\"],\n  [\"25+52=77\", \"47+51=98\"],\n  [\"62+8=70\", \"34+55=89\"],\n  [\"43+16=59\", \"9+83=92\"],\n  [\"87-82=5\", \"21-13=8\"],\n  [\"93-79=14\", \"76-7=69\"],\n  [\"42+36=78\", \"83-34=49\"],\n  [\"13+36=49\", \"85-42=43\"],\n  [\"33-31=2\", \"49-16=33\"],\n  [\"10+36=46\", \"96-89=7\"],\n  [\"6+67=73\", \"23-10=13\"],\n  [\"89-67=22\", \"31+8=39\"],\n  [\"36-35=1\", \"95-24=71\"],\n  [\"12-12=0\", \"62-17=45\"],\n  [\"90-72=18\", \"36+30=66\"],\n  [\"11+43=54\", \"38+25=63\"],\n  [\"82+3=85\", \"11+63=74\"],\n  [\"63-63=0\", \"40+44=84\"],\n  [\"12-5=7\", \"74-31=43\"],\n  [\"38+4=42\", \"84+14=98\"],\n  [\"2+5=7\", \"95-8=87\"],\n  [\"64-17=47\", \"87-68=19\"],\n  [\"53-13=40\", \"73-30=43\"],\n  [\"67+25=92\", \"37+36=73\"],\n  [\"71+7=78\", \"41-17=24\"],\n  [\"2+89=91\", \"42-30=12\"],\n  [\"55-45=10\", \"36+44=80\"],\n  [\"15-6=9\", \"34-32=2\"],\n  [\"17+71=88\", \"76+15=91\"],\n  [\"84-79=5\", \"63-56=7\"],\n  [\"56+5=61\", \"94-37=57\"],\n  [\"27+10=37\", \"18+24=42\"],\n  [\"42-25=17\", \"47-1=46\"],\n  [\"5-3=2\", \"90-55=35\"],\n  [\"45-1=44\", \"31+54=85\"],\n  [\"37+49=86\", \"81-32=49\"],\n  [\"66+24=90\", \"12-10=2\"],\n  [\"1+88=89\", \"42+53=95\"],\n  [\"12+15=27\", \"81-75=6\"],\n  [\"89+8=97\", \"34+43=77\"],\n  [\"44+5=49\", \"1+98=99\"],\n  [\"65+5=70\", \"41+4=45\"],\n  [\"18+27=45\", \"86-43=43\"],\n  [\"7+76=83\", \"29+2=31\"],\n  [\"93-91=2\", \"33-29=4\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error('Expected exactly 1 match for \"' + oldText + '\" but found ' + results.items.length);\n  }\n  results.items[0].insertText(newText, 'Replace');\n}\nawait context.sync();", "ps1": "# Update the header date and every math-fact answer in the table.\n$d = $word.ActiveDocument\n\n# --- Header date paragraph (first paragraph of the body) ---\n$dateOld = \"2023-11-23 Thursday\"\n$dateNew = \"2023-11-24 Friday\"\n$dateRange = $d.Paragraphs.Item(1).Range\nif ($dateRange.Text.TrimEnd([char]13, [char]7) -ne $dateOld) {\n    throw \"Unexpected header text: '$($dateRange.Text)'\"\n}\n$dateRange.Text = $dateNew\n\n# --- Table of math facts: 20 rows x 5 columns ---\n$t = $d.Tables.Item(1)\n\n$grid = @(\n    @(@(\"64+2=66\", \"15+64=79\"), @(\"11+47=58\", \"48-5=43\"), @(\"22+53=75\", \"40-27=13\"), @(\"91-32=59\", \"94-27=67\"), @(\"32+20=52\", \"94-42=52\")),\n    @(@(\"69-11=58\", \"15+32=47\"), @(\"62-61=1\", \"17+79=96\"), @(\"18+34=52\", \"43+14=57\"), @(\"4+26=30\", \"85+1=86\"), @(\"43+11=54\", \"59-59=0\")),\n    @(@(\"55+39=94\", \"78-13=65\"), @(\"10-10=0\", \"53-4=49\"), @(\"50+0=50\", \"30+7=37\"), @(\"79+10=89\", \"21+32=53\"), @(\"28-24=4\", \"65+9=74\")),\n    @(@(\"56-27=29\", \"49+13=62\"), @(\"54-50=4\", \"91-29=62\"), @(\"84-24=60\", \"63-45=18\"), @(\"38+5=43\", \"6+7=13\"), @(\"22+22=44\", \"57+21=78\")),\n    @(@(\"59+23=82\", \"80-44=36\"), @(\"64-12=52\", \"39-22=17\"), @(\"8+71=79\", \"22+38=60\"), @(\"28-15=13\", \"23+8=31\"), @(\"34+21=55\", \"4+75=79\")),\n    @(@(\"3+61=64\", \"67-22=45\"), @(\"14+4=18\", \"94-73=21\"), @(\"87-35=52\", \"63+15=78\"), @(\"28+40=68\", \"71+1=72\"), @(\"63+32=95\", \"10+13=23\")),\n    @(@(\"88-26=62\", \"78-6=72\"), @(\"79-46=33\", \"71-17=54\"), @(\"75-75=0\", \"7+12=19\"), @(\"99-6=93\", \"57-33=24\"), @(\"60-12=48\", \"54-8=46\")),\n    @(@(\"27+36=63\", \"64+28=92\"), @(\"28+24=52\", \"18+52=70\"), @(\"5+76=81\", \"25+63=88\"), @(\"45+31=76\", \"68-3=65\"), @(\"27+58=85\", \"94-29=65\")),\n    @(@(\"61+7=68\", \"20-1=19\"), @(\"98-65=33\", \"61+34=95\"), @(\"44-35=9\", \"34+59=93\"), @(\"51-2=49\", \"70-52=18\"), @(\"4+34=38\", \"65-46=19\")),\n    @(@(\"94+3=97\", \"74-34=40\"), @(\"23+73=96\", \"69-30=39\"), @(\"48-27=21\", \"27+56=83\"), @(\"25+51=76\", \"29+19=48\"), @(\"83-8=75\", \"39+42=81\")),\n    @(@(\"94-81=13\", \"67-10=57\"), @(\"9+42=51\", \"2+26=28\"), @(\"49-38=11\", \"58+20=78\"), @(\"81-23=58\", \"46+16=62\"), @(\"9+28=37\", \"93-45=48\")),\n    @(@(\"34-12=22\", \"88-47=41\"), @(\"25+52=77\", \"47+51=98\"), @(\"62+8=70\", \"34+55=89\"), @(\"43+16=59\", \"9+83=92\"), @(\"87-82=5\", \"21-13=8\")),\n    @(@(\"93-79=14\", \"76-7=69\"), @(\"42+36=78\", \"83-34=49\"), @(\"13+36=49\", \"85-42=43\"), @(\"33-31=2\", \"49-16=33\"), @(\"10+36=46\", \"96-89=7\")),\n    @(@(\"6+67=73\", \"23-10=13\"), @(\"89-67=22\", \"31+8=39\"), @(\"36-35=1\", \"95-24=71\"), @(\"12-12=0\", \"62-17=45\"), @(\"90-72=18\", \"36+30=66\")),\n    @(@(\"11+43=54\", \"38+25=63\"), @(\"82+3=85\", \"11+63=74\"), @(\"63-63=0\", \"40+44=84\"), @(\"12-5=7\", \"74-31=43\"), @(\"38+4=42\", \"84+14=98\")),\n    @(@(\"2+5=7\", \"95-8=87\"), @(\"64-17=47\", \"87-68=19\"), @(\"53-13=40\", \"73-30=43\"), @(\"67+25=92\", \"37+36=73\"), @(\"71+7=78\", \"41-17=24\")),\n    @(@(\"2+89=91\", \"42-30=12\"), @(\"55-45=10\", \"36+44=80\"), @(\"15-6=9\", \"34-32=2\"), @(\"17+71=88\", \"76+15=91\"), @(\"84-79=5\", \"63-56=7\")),\n    @(@(\"56+5=61\", \"94-37=57\"), @(\"27+10=37\", \"18+24=42\"), @(\"42-25=17\", \"47-1=46\"), @(\"5-3=2\", \"90-55=35\"), @(\"45-1=44\", \"31+54=85\")),\n    @(@(\"37+49=86\", \"81-32=49\"), @(\"66+24=90\", \"12-10=2\"), @(\"1+88=89\", \"42+53=95\"), @(\"12+15=27\", \"81-75=6\"), @(\"89+8=97\", \"34+43=77\")),\n    @(@(\"44+5=49\", \"1+98=99\"), @(\"65+5=70\", \"41+4=45\"), @(\"18+27=45\", \"86-43=43\"), @(\"7+76=83\", \"29+2=31\"), @(\"93-91=2\", \"33-29=4\")),\n)\n\nfor ($r = 0; $r -lt $grid.Count; $r++) {\n    $row = $grid[$r]\n    for ($c = 0; $c -lt $row.Count; $c++) {\n        $old = $row[$c][0]\n        $new = $row[$c][1]\n        $cell = $t.Cell($r + 1, $c + 1)\n        $cellRange = $cell.Range\n        $actual = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($actual -ne $old) {\n            throw \"Row $($r+1) Col $($c+1): expected '$old' but found '$actual'\"\n        }\n        $cellRange.Text = $new\n    }\n}\n\nWrite-Output 'Done'"}
